$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds the new values for columns B,C,D,E,F,I,J,K,L,M,N
# for rows 2 through 25 (in order), matching the case with 380 kV results.
$data = @(
    @(1.02,1.012292433325808,1.019248255289404,0.9926147277508489,1.010565194270536,1.025142902596617,1.017535986546285,1.022092065907909,0.9955398523336033,1.013435001458975,1.010006368440002),
    @(1.02,1.013194884753648,1.019912264129034,0.9936372048519304,1.01211306238351,1.025279182911297,1.018072321354631,1.022562449592036,0.9963617723202692,1.014784889401855,1.010182616354473),
    @(1.02,1.013778205724422,1.020340805088882,0.9942998659930995,1.013113891899843,1.025365041537151,1.018418157468102,1.022865026159951,0.9968940712668345,1.015657145157266,1.010296256525005),
    @(1.02,1.014023284359751,1.020520695589501,0.9945786998346017,1.013534468210394,1.025400579705842,1.018563257768412,1.022991799690596,0.997117960005301,1.016023556294678,1.010343934143808),
    @(1.02,1.01406442538067,1.020550884253621,0.9946255319796338,1.013605074882438,1.025406514056117,1.018587603790518,1.023013060317406,0.9971555583673453,1.016085061757761,1.010351933755021),
    @(1.02,1.013781481064326,1.020343209848974,0.9943035907982488,1.013119512331341,1.02536551858835,1.018420097441717,1.022866721802865,0.9968970624462087,1.015662042275254,1.010296893975678),
    @(1.02,1.012597550660024,1.019472891106258,0.9929600610674301,1.011088461264505,1.025189440061672,1.017717493663473,1.022251405333723,0.995817528259106,1.013891456366409,1.010066015911707),
    @(1.02,1.010506511959808,1.017930750186949,0.9906006454969559,1.007503523466787,1.0248613888231,1.01647015626887,1.02115342122297,0.9939188001724441,1.010761951395903,1.009656083963866),
    @(1.02,1.009109236636048,1.016896956600626,0.989033133672735,1.005109183150741,1.024630759183526,1.015632361699444,1.020412225950081,0.9926553831429383,1.008668907597992,1.009380712308816),
    @(1.02,1.008503424292771,1.016447966764536,0.988355674866747,1.004071284923187,1.024528069152499,1.015268104828417,1.02008909925757,0.9921088820399291,1.007760934668866,1.009260978190957),
    @(1.02,1.008278280339407,1.016280989434228,0.9881042295826724,1.003685585824752,1.024489501224578,1.015132579900423,1.019968747536939,0.9919059725120875,1.007423415805989,1.009216428915346),
    @(1.02,1.008326579844561,1.016316815796937,0.9881581567098651,1.003768327739182,1.02449779337121,1.015161660587524,1.019994578234822,0.9919494934313052,1.00749582649101,1.009225988268176),
    @(1.02,1.008484816239976,1.016434168495512,0.9883348863814464,1.004039406578152,1.024524889769138,1.01525690686159,1.020079157639151,0.9920921077337197,1.007733040539127,1.009257297258361),
    @(1.02,1.008582295213177,1.016506446541734,0.9884438009545853,1.004206403461338,1.024541528541282,1.015315561578596,1.020131226329853,0.9921799884222134,1.007879161663091,1.009276577859333),
    @(1.02,1.009149425855894,1.016926726158035,0.9890781214508737,1.005178040547414,1.024637514868774,1.01565650487936,1.020433624792093,0.9926916645766087,1.008729131032007,1.00938864819988),
    @(1.02,1.009504961864336,1.017189995186735,0.989476357848556,1.005787214363248,1.024696968189819,1.015869971398944,1.020622726635695,0.9930127773699352,1.009261842572634,1.009458813988622),
    @(1.02,1.00971226450854,1.017343425240594,0.9897087662937556,1.006142426307959,1.024731373632327,1.015994339430905,1.020732815789373,0.9932001317071769,1.009572403378979,1.00949969262148),
    @(1.02,1.009782936567418,1.017395718808509,0.9897880325774034,1.006263526161091,1.024743058740236,1.016036721408711,1.02077031759584,0.9932640239640975,1.009678269498521,1.009513623060337),
    @(1.02,1.009466824028588,1.017161762387903,0.9894336180360679,1.005721867086884,1.024690617612947,1.015847083294567,1.020602459595867,0.9929783193494215,1.00920470434463,1.009451290817254),
    @(1.02,1.008438222835778,1.016399616638268,0.9882828385668249,1.003959585574474,1.024516922265303,1.015228865386626,1.020054260168746,0.9920501090198102,1.007663194066968,1.009248079603984),
    @(1.02,1.007790815799254,1.015919253837133,0.9875604150241495,1.002850538745415,1.024405258490922,1.014838872732139,1.019707687285607,0.9914670000341481,1.006692493663826,1.009119880434558),
    @(1.02,1.008134083637929,1.016174014177117,0.9879432794643023,1.003438565568697,1.024464686150774,1.015045738064096,1.019891592030072,0.991776070289318,1.007207223799417,1.009187882250011),
    @(1.02,1.009484057107896,1.017174519976965,0.9894529299347244,1.00575139504279,1.024693488007733,1.015857425886511,1.020611618049815,0.9929938892766442,1.009230523145781,1.009454690362447),
    @(1.02,1.011047666730574,1.018330437589426,0.9912096547607049,1.00843106083834,1.024948302193921,1.016793721230952,1.02143890021102,0.9944092447426414,1.009762428250256,1.009296893975678)
)

# Columns in the order the values appear within each row above.
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $rowVals[$j]
    }
}
